$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "will tomorrow rise" column (J) originally compared tomorrow's Close (G)
# to today's Close (G). The formula is changed so it instead compares
# tomorrow's Close (G) to tomorrow's Open (D).
#
# Original: =IF(SIGN(G3-G2)=1,1,0)   (in J2)
# New:      =IF(SIGN(G3-D3)=1,1,0)   (in J2)
#
# Apply this across the whole data range (J2:J2077), letting Excel adjust
# the relative references for each row, which mirrors the shared-formula
# behaviour seen in the workbook.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # -4162 = xlUp, column G = 7

$ws.Range("J2:J$lastRow").Formula = "=IF(SIGN(G3-D3)=1,1,0)"
